# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets
# to reflect the newly scraped counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new value for column F
$updates = @{
    4  = 10336
    8  = 7136
    13 = 3185
    18 = 1034
    19 = 276
    20 = 65
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
